$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46026, 2),
    @(46021, 2),
    @(46022, 10),
    @(46025, 2),
    @(46027, 2),
    @(46020, 4),
    @(46015, 7),
    @(46018, 6),
    @(46016, 4),
    @(46024, 10)
)

$startRow = 43
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $serial = $data[$i][0]
    $count = $data[$i][1]

    # Copy the formatting (date number format, style) from an existing
    # date cell in column A so the new rows match the existing ones.
    $ws.Cells.Item(2, 1).Copy($ws.Cells.Item($row, 1))
    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = $count
}

$ws.Range("A43:B52").Select()
